$d = $word.ActiveDocument

# 1) Remove the stray "_GoBack" bookmark that sits alone in its own
#    (otherwise empty) paragraph, right after the big blue rectangle
#    shape.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) In the "Discretized version of the equations <bullet> Description
#    of the numerical method" heading paragraph, the " • " bullet
#    separator becomes a plain " " space.
$d.Content.Find.Execute(" " + [char]0x2022 + " ", $true, $false, $false, `
    $false, $false, $true, 1, $false, " ", 2) | Out-Null

# 3) Re-insert the "_GoBack" bookmark right after that (now plain)
#    space run, immediately before "Description of the numerical
#    method".
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.IndexOf("Discretized") -ge 0) {
        $target = $p
        break
    }
}

$relIdx = $target.Range.Text.IndexOf("Description")
$pos = $target.Range.Start + $relIdx
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
